# Edit script for NYPD CompStat 120th Precinct weekly report
# Updates the report header (volume/date range) and refreshes the crime-stat table
# for the new reporting week (new weekly crime data collected).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (sharedStrings rich-text cells) ---
$ws.Range("A8").Value = "Volume 31   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/15/2024  Through  1/21/2024"

# --- Reference cells used as format donors (stable styles already present in sheet) ---
$refText0 = $ws.Range("C14")      # style: right-aligned General text "0"
$refTextStars = $ws.Range("M14")  # style: right-aligned General text "***.*"
$refNum15 = $ws.Range("G14")      # style: #,##0 integer
$refNum16 = $ws.Range("H14")      # style: #,##0.0 percentage-change number

# --- Cells that change from a number to the literal placeholder text ("0" or "***.*") ---
$ws.Range("D14").Value = "'0"
$refText0.Copy()
$ws.Range("D14").PasteSpecial(-4122)
$ws.Range("E14").Value = "'***.*"
$refTextStars.Copy()
$ws.Range("E14").PasteSpecial(-4122)
$ws.Range("C15").Value = "'0"
$refText0.Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C23").Value = "'0"
$refText0.Copy()
$ws.Range("C23").PasteSpecial(-4122)

# --- Cells that change from placeholder text to a real number ---
$ws.Range("L14").Value = -100
$refNum16.Copy()
$ws.Range("L14").PasteSpecial(-4122)
$ws.Range("D15").Value = 1
$refNum15.Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100
$refNum16.Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("G15").Value = 1
$refNum15.Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("H15").Value = 0
$refNum16.Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("J15").Value = 1
$refNum15.Copy()
$ws.Range("J15").PasteSpecial(-4122)
$ws.Range("K15").Value = 0
$refNum16.Copy()
$ws.Range("K15").PasteSpecial(-4122)
$ws.Range("M23").Value = 200
$refNum16.Copy()
$ws.Range("M23").PasteSpecial(-4122)
$ws.Range("D26").Value = 1
$refNum15.Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("E26").Value = 0
$refNum16.Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("G26").Value = 1
$refNum15.Copy()
$ws.Range("G26").PasteSpecial(-4122)
$ws.Range("H26").Value = 300
$refNum16.Copy()
$ws.Range("H26").PasteSpecial(-4122)
$ws.Range("J26").Value = 1
$refNum15.Copy()
$ws.Range("J26").PasteSpecial(-4122)
$ws.Range("K26").Value = 100
$refNum16.Copy()
$ws.Range("K26").PasteSpecial(-4122)
$ws.Range("L28").Value = 0
$refNum16.Copy()
$ws.Range("L28").PasteSpecial(-4122)
$ws.Range("L29").Value = 0
$refNum16.Copy()
$ws.Range("L29").PasteSpecial(-4122)

# --- Plain numeric updates (style unchanged) ---
$ws.Range("F15").Value = 1
$ws.Range("N15").Value = -87.5
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 6
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 7
$ws.Range("G16").Value = 23
$ws.Range("H16").Value = -69.565217391304
$ws.Range("I16").Value = 5
$ws.Range("J16").Value = 20
$ws.Range("K16").Value = -75
$ws.Range("L16").Value = -28.571428571428
$ws.Range("M16").Value = -66.666666666666
$ws.Range("N16").Value = -93.243243243243
$ws.Range("C17").Value = 4
$ws.Range("D17").Value = 12
$ws.Range("E17").Value = -66.666666666666
$ws.Range("F17").Value = 35
$ws.Range("G17").Value = 34
$ws.Range("H17").Value = 2.941176470588
$ws.Range("I17").Value = 20
$ws.Range("J17").Value = 25
$ws.Range("K17").Value = -20
$ws.Range("L17").Value = 5.263157894736
$ws.Range("M17").Value = 42.857142857142
$ws.Range("N17").Value = -41.176470588235
$ws.Range("C18").Value = 2
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 8
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = -42.857142857142
$ws.Range("I18").Value = 5
$ws.Range("J18").Value = 10
$ws.Range("L18").Value = 66.666666666666
$ws.Range("M18").Value = -68.75
$ws.Range("N18").Value = -94.444444444444
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 11
$ws.Range("E19").Value = -18.181818181818
$ws.Range("F19").Value = 24
$ws.Range("G19").Value = 35
$ws.Range("H19").Value = -31.428571428571
$ws.Range("I19").Value = 16
$ws.Range("J19").Value = 24
$ws.Range("K19").Value = -33.333333333333
$ws.Range("L19").Value = -48.387096774193
$ws.Range("M19").Value = 6.666666666666
$ws.Range("N19").Value = -44.827586206896
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 7
$ws.Range("H20").Value = 16.666666666666
$ws.Range("I20").Value = 6
$ws.Range("J20").Value = 4
$ws.Range("K20").Value = 50
$ws.Range("L20").Value = 500
$ws.Range("M20").Value = -25
$ws.Range("N20").Value = -93.181818181818
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 36
$ws.Range("E21").Value = -44.444444444444
$ws.Range("F21").Value = 82
$ws.Range("G21").Value = 114
$ws.Range("H21").Value = -28.070175438596
$ws.Range("I21").Value = 53
$ws.Range("J21").Value = 85
$ws.Range("K21").Value = -37.647058823529
$ws.Range("L21").Value = -15.873015873015
$ws.Range("M21").Value = -23.188405797101
$ws.Range("N21").Value = -83.641975308642
$ws.Range("D23").Value = 4
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 6
$ws.Range("G23").Value = 13
$ws.Range("H23").Value = -53.846153846153
$ws.Range("I23").Value = 3
$ws.Range("J23").Value = 9
$ws.Range("K23").Value = -66.666666666666
$ws.Range("L23").Value = -25
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 31
$ws.Range("E24").Value = -25.806451612903
$ws.Range("F24").Value = 99
$ws.Range("G24").Value = 104
$ws.Range("H24").Value = -4.807692307692
$ws.Range("I24").Value = 74
$ws.Range("J24").Value = 80
$ws.Range("K24").Value = -7.5
$ws.Range("L24").Value = 51.020408163265
$ws.Range("M24").Value = -3.896103896103
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = -50
$ws.Range("F25").Value = 41
$ws.Range("G25").Value = 49
$ws.Range("H25").Value = -16.326530612244
$ws.Range("I25").Value = 33
$ws.Range("J25").Value = 40
$ws.Range("K25").Value = -17.5
$ws.Range("L25").Value = 10
$ws.Range("M25").Value = -50.746268656716
$ws.Range("I26").Value = 2
$ws.Range("L26").Value = 0
$ws.Range("C27").Value = 3
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 200
$ws.Range("F27").Value = 5
$ws.Range("H27").Value = 25
$ws.Range("I27").Value = 5
$ws.Range("J27").Value = 4
$ws.Range("K27").Value = 25
$ws.Range("L27").Value = 66.666666666666
$ws.Range("M28").Value = -66.666666666666
$ws.Range("M29").Value = -66.666666666666

$excel.CutCopyMode = $false
